# Auto-generated Excel COM-interop script
# Applies numeric updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# as described by the upstream OOXML diff (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 87
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 24
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 24
$ws.Range("M5").Value = -35
$ws.Range("N5").Value = -254

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H137").Value = 2430.1667
$ws.Range("J137").Value = 3212.1904
$ws.Range("L137").Value = 9636.5712
$ws.Range("N137").Value = -14736.5712

$ws.Range("H138").Value = 2997.25
$ws.Range("J138").Value = 2994
$ws.Range("L138").Value = 8982
$ws.Range("N138").Value = -19262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 670.5
$ws.Range("I5").Value = 739.8
$ws.Range("J5").Value = 555
$ws.Range("K5").Value = 739.8
$ws.Range("L5").Value = 555
$ws.Range("M5").Value = -627.8
$ws.Range("N5").Value = -779

$ws.Range("H16").Value = 590
$ws.Range("I16").Value = 590
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 590
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -303
$ws.Range("N16").ClearContents()

$ws.Range("H74").Value = 3143.2273
$ws.Range("I74").Value = 2169
$ws.Range("K74").Value = 2169
$ws.Range("M74").Value = -1295

$ws.Range("H77").Value = 3143.2273
$ws.Range("I77").Value = 2169
$ws.Range("K77").Value = 10845
$ws.Range("M77").Value = -6477

$ws.Range("H110").Value = 1151.25
$ws.Range("I110").Value = 1173
$ws.Range("J110").Value = 999
$ws.Range("K110").Value = 1173
$ws.Range("L110").Value = 999
$ws.Range("M110").Value = 872
$ws.Range("N110").Value = -5089

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 670.5
$ws.Range("I4").Value = 739.8
$ws.Range("J4").Value = 555
$ws.Range("K4").Value = 739.8
$ws.Range("L4").Value = 555
$ws.Range("M4").Value = -624.8
$ws.Range("N4").Value = -785

$ws.Range("H64").Value = 1332
$ws.Range("J64").Value = 1664
$ws.Range("L64").Value = 1664
$ws.Range("N64").Value = -2114

$ws.Range("H67").Value = 1332
$ws.Range("J67").Value = 1664
$ws.Range("L67").Value = 1664
$ws.Range("N67").Value = -3224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 774.75
$ws.Range("I22").Value = 774.75
$ws.Range("K22").Value = 774.75
$ws.Range("M22").Value = -424.75

$ws.Range("H58").Value = 2883.2
$ws.Range("J58").Value = 2499
$ws.Range("L58").Value = 2499
$ws.Range("N58").Value = -2905

$ws.Range("H134").Value = 2324.0557
$ws.Range("I134").Value = 1988.9333
$ws.Range("K134").Value = 5966.7999
$ws.Range("M134").Value = -3431.7999

$ws.Range("H136").Value = 2883.2
$ws.Range("J136").Value = 2499
$ws.Range("L136").Value = 7497
$ws.Range("N136").Value = -12597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 840.6667
$ws.Range("I5").Value = 911
$ws.Range("K5").Value = 2733
$ws.Range("M5").Value = -2621

$ws.Range("H8").Value = 460
$ws.Range("I8").Value = 460
$ws.Range("K8").Value = 1380
$ws.Range("M8").Value = -1241

$ws.Range("H23").Value = 281
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 323.66666
$ws.Range("K23").Value = 75
$ws.Range("L23").Value = 970.9999799999999
$ws.Range("M23").Value = 160
$ws.Range("N23").Value = -1440.99998

$ws.Range("H131").Value = 2170.875
$ws.Range("I131").Value = 1499
$ws.Range("J131").Value = 2266.8572
$ws.Range("K131").Value = 4497
$ws.Range("L131").Value = 6800.571599999999
$ws.Range("M131").Value = 543
$ws.Range("N131").Value = -16880.5716

$ws.Range("H132").Value = 4860.9414
$ws.Range("J132").Value = 5365.231
$ws.Range("L132").Value = 48287.079
$ws.Range("N132").Value = -53347.079

$ws.Range("H134").Value = 1577.2
$ws.Range("I134").Value = 1577.2
$ws.Range("K134").Value = 4731.6
$ws.Range("M134").Value = 338.3999999999996

$ws.Range("H135").Value = 840.6667
$ws.Range("I135").Value = 911
$ws.Range("K135").Value = 8199
$ws.Range("M135").Value = -5664

$ws.Range("H137").Value = 5784.857
$ws.Range("J137").Value = 5784.857
$ws.Range("L137").Value = 17354.571
$ws.Range("N137").Value = -27554.571

$ws.Range("H138").Value = 5074.25
$ws.Range("I138").Value = 5074.25
$ws.Range("K138").Value = 15222.75
$ws.Range("M138").Value = -10082.75

$ws.Range("H139").Value = 1987.4
$ws.Range("I139").Value = 1987.4
$ws.Range("K139").Value = 5962.200000000001
$ws.Range("M139").Value = -822.2000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4381.6665
$ws.Range("I61").Value = 4358
$ws.Range("K61").Value = 4358
$ws.Range("M61").Value = -4156

$ws.Range("H113").Value = 4381.6665
$ws.Range("I113").Value = 4358
$ws.Range("K113").Value = 4358
$ws.Range("M113").Value = -2188

$ws.Range("H136").Value = 6591.5293
$ws.Range("I136").Value = 6566
$ws.Range("K136").Value = 19698
$ws.Range("M136").Value = -17148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25000
$ws.Range("J2").Value = 25000
$ws.Range("L2").Value = 25000
$ws.Range("N2").Value = -25224

$ws.Range("H122").Value = 2912.5386
$ws.Range("I122").Value = 2864.1667
$ws.Range("K122").Value = 8592.500100000001
$ws.Range("M122").Value = -6142.500100000001

$ws.Range("H126").Value = 2831.8333
$ws.Range("I126").Value = 2831.8333
$ws.Range("K126").Value = 8495.499899999999
$ws.Range("M126").Value = -6025.499899999999

$ws.Range("H132").Value = 3608.9
$ws.Range("I132").Value = 3348
$ws.Range("K132").Value = 10044
$ws.Range("M132").Value = -7514

$ws.Range("H136").Value = 1384.5454
$ws.Range("I136").Value = 1384.5454
$ws.Range("K136").Value = 4153.6362
$ws.Range("M136").Value = -1603.6362
